$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values formatted as text (e.g. "1.020", "27.682.15")
# Setting .Value on a numeric-looking string would make Excel silently
# coerce it to a real number (dropping trailing zeros / using sci notation),
# so every touched D cell is forced to Text format first.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.682.15"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.859.36"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  -0.93%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.89"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4367"
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3775"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07412"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8828"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.56"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.862.57"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.749"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.489"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07144"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.47"
$ws.Range("E16").Value = "  +5.29%  "
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009043"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.019"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.46"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.696.79"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.293"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.15"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.096.01"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.037"
$ws.Range("E25").Value = "  +6.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.05"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.70"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.414"
$ws.Range("E28").Value = "  +2.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.992"
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.69"
$ws.Range("E30").Value = "  +3.96%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.214"
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7681"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("E34").Value = "  +5.28%  "
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.140"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01974"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05305"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.870"
$ws.Range("E40").Value = "  +1.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5177"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.976"
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1678"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.712"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.78"
$ws.Range("E45").Value = "  +1.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "110.11"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.711"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4727"
$ws.Range("E48").Value = "  +1.80%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.020"
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06479"
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.839"
$ws.Range("E51").Value = "  -0.45%  "
